# Zowe HA Installation Roadmap - WIP edit
# Split the "(One-time setup per z/OS environment) " run on slide 3's
# "Rectangle 86" shape so that "z/OS" becomes "Sysplex", matching the
# rest of the deck's existing "(One-time setup per Sysplex environment)"
# wording elsewhere on the same slide.

$p = $ppt.ActivePresentation

# The shape lives on the 3rd slide (ppt/slides/slide3.xml).
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item("Rectangle 86")

$tr = $shape.TextFrame.TextRange

# Current text begins: "(One-time setup per z/OS environment) Create the ..."
# Replace the "z/OS " substring (including its trailing space) with
# "Sysplex " so the text reads "...setup per Sysplex environment)...".
# Characters(Start, Length) is 1-based; "z/OS " starts right after
# "(One-time setup per " (21 characters in), and is 5 characters long.
$target = $tr.Characters(21, 5)
$target.Text = "Sysplex "
